$d = $word.ActiveDocument

# The document currently ends with a single empty paragraph (right before the
# sectPr). We replace that empty paragraph with two new paragraphs:
#   1) A bold, 16pt ("sz"=32 half-points) heading paragraph reading "Scrum:"
#   2) A body paragraph describing the decision to use scrum, with the word
#      "scrum" marked (twice) with spell-check proofErr wrappers, matching
#      the exact run/proofErr layout produced by Word's own spell checker.
$target = $d.Paragraphs.Item($d.Paragraphs.Count)
$targetRange = $target.Range

$w = 'http://schemas.openxmlformats.org/wordprocessingml/2006/main'

$headingXml = '<w:p xmlns:w="' + $w + '">' + `
  '<w:pPr><w:rPr><w:b/><w:bCs/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr></w:pPr>' + `
  '<w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="32"/><w:szCs w:val="32"/></w:rPr><w:t>Scrum:</w:t></w:r>' + `
  '</w:p>'

$bodyXml = '<w:p xmlns:w="' + $w + '">' + `
  '<w:r><w:t xml:space="preserve">Eu decidi usar </w:t></w:r>' + `
  '<w:proofErr w:type="spellStart"/>' + `
  '<w:r><w:t>scrum</w:t></w:r>' + `
  '<w:proofErr w:type="spellEnd"/>' + `
  '<w:r><w:t xml:space="preserve"> no meu projeto por ser um dos métodos de framework mais ágeis no mercado onde ele visa sempre a eficiência do projeto e a coerência com o que o dono do negocio quer, no método </w:t></w:r>' + `
  '<w:proofErr w:type="spellStart"/>' + `
  '<w:r><w:t>scrum</w:t></w:r>' + `
  '<w:proofErr w:type="spellEnd"/>' + `
  '<w:r><w:t xml:space="preserve"> estamos de contato direto com o dono do </w:t></w:r>' + `
  '<w:r><w:t>negócio</w:t></w:r>' + `
  '<w:r><w:t xml:space="preserve"> assim decidindo as prioridades do projeto</w:t></w:r>' + `
  '<w:r><w:t xml:space="preserve"> e a ordem do que será desenvolvido primeiro</w:t></w:r>' + `
  '<w:r><w:t xml:space="preserve"> juntamente a quem vai usa-lo garantindo assim um numero de falhas muito baixo.</w:t></w:r>' + `
  '</w:p>'

$null = $targetRange.InsertXML($headingXml + $bodyXml)

Write-Output "Inserted Scrum heading + body paragraphs"
